# Rename roles -> capabilities, role -> capability
#
# The "Database" sheet documents a `roles` table (and the `user_role`
# pivot table that links users to roles). This rename turns that into a
# `capabilities` table / `user_capability` pivot table, updating every
# cell that mentions "role"/"roles" to "capability"/"capabilities".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

# Table name cell (row 12) and its use as a FK target description (row 18)
$ws.Range("A12").Value = "capabilities"
$ws.Range("D18").Value = "capabilities"

# Long-form description of the roles/capabilities table (row 13)
$ws.Range("A13").Value = "Danh sách các quyền theo" + [char]10 + "thứ tự kèm theo tên. Có" + [char]10 + "23 capabilities được định " + [char]10 + "nghĩa trong types.php"

# user_role -> user_capability pivot table name + its FK column (row 14)
$ws.Range("A14").Value = "user_capability"
$ws.Range("C14").Value = "capability_id"

# Description of how the pivot table works (row 15)
$ws.Range("A15").Value = "Nếu user A (id=2) có capability B (id=5) thì cặp" + [char]10 + " (3,5) được lưu vào đây" + [char]10 + "(Đọc types.php) để hiểu kỹ"

# FK note on the pivot table's capability_id column (row 15, col C)
$ws.Range("C15").Value = "FK:capabilities.id"

# Match the author's final selection/scroll position in the saved view
$ws.Range("C15").Select()
